# 5th commit by Pradnya
# - Re-ran the AWB test suite: updated the execution timestamp and the
#   AWB reference numbers embedded in a couple of the "Actual Result1"
#   messages, and test_View_AWB_TC1 (row 6) now fails with a new assertion
#   message instead of passing.
# - Widened the Execution Date1 / Actual Result1 / Status1 columns so the
#   longer messages are readable, and moved the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates -------------------------------------------------

# Execution Date1 column (E2:E8 all share the same value in this sheet)
$ws.Range("E2:E8").Value = "21-Nov-2025 03:42:59 PM"

# Actual Result1 text for the PP AWB / CC AWB test cases (AWB numbers changed)
$ws.Range("F2").Value = "AWB 057 97024885 did not validate successfully."
$ws.Range("F3").Value = "AWB 057 97024896 did not validate successfully."

# test_View_AWB_TC1 (row 6) now fails, with a captured actual-result message
$ws.Range("F6").Value = "Expected '0 records found' but got ' 1 records found'"
$ws.Range("G6").Value = "Fail"

# --- Column width adjustments ----------------------------------------------
# (ColumnWidth is in character units and gets snapped to the host's pixel
#  grid, so the inputs below are chosen to land on the exact target OOXML
#  widths of 22.36328125 / 65.08984375 / 10.90625 after that snap.)

$ws.Columns.Item(5).ColumnWidth = 21.5                # column E
$ws.Columns.Item(6).ColumnWidth = 64.333333333333336  # column F
$ws.Columns.Item(7).ColumnWidth = 10                  # column G

# --- Selection / view changes ----------------------------------------------

$ws.Range("F6").Select() | Out-Null
